# Refresh the cryptos worksheet (price / 1h volume columns, plus a few
# reordered coin rows) with the latest scraped values from the feed,
# matching the automated "Updated cryptos list ... with GitHub Actions" run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even when it looks numeric
# (e.g. '8.90', '1.00', '3.453.58'), so formatting/precision such as
# trailing zeros or the site's "."-grouped price strings survive --
# exactly like the source inline-string cells already stored them.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '69.622.24'
Set-TextValue $ws.Range('E2') '  -0.44%  '

Set-TextValue $ws.Range('D3') '3.453.58'
Set-TextValue $ws.Range('E3') '  -1.68%  '

Set-TextValue $ws.Range('D4') '0.998'
Set-TextValue $ws.Range('E4') '  -0.16%  '

Set-TextValue $ws.Range('D5') '610.62'
Set-TextValue $ws.Range('E5') '  +1.24%  '

Set-TextValue $ws.Range('D6') '167.06'
Set-TextValue $ws.Range('E6') '  -2.95%  '

Set-TextValue $ws.Range('D7') '3.447.73'
Set-TextValue $ws.Range('E7') '  -1.69%  '

Set-TextValue $ws.Range('D8') '0.594'
Set-TextValue $ws.Range('E8') '  -2.33%  '

Set-TextValue $ws.Range('E9') '  +0.08%  '

Set-TextValue $ws.Range('E10') '  +0.05%  '

Set-TextValue $ws.Range('D11') '7.09'
Set-TextValue $ws.Range('E11') '  -2.23%  '

Set-TextValue $ws.Range('D12') '0.563'
Set-TextValue $ws.Range('E12') '  -2.94%  '

Set-TextValue $ws.Range('D13') '44.35'
Set-TextValue $ws.Range('E13') '  -3.32%  '

Set-TextValue $ws.Range('D14') '0.0000268'
Set-TextValue $ws.Range('E14') '  -2.11%  '

Set-TextValue $ws.Range('D15') '4.013.43'
Set-TextValue $ws.Range('E15') '  -1.71%  '

Set-TextValue $ws.Range('D16') '8.18'
Set-TextValue $ws.Range('E16') '  -0.81%  '

Set-TextValue $ws.Range('B17') 'BitcoinCash'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D17') '584.48'
Set-TextValue $ws.Range('E17') '  -3.20%  '

Set-TextValue $ws.Range('B18') 'WrappedEther'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D18') '3.458.47'
Set-TextValue $ws.Range('E18') '  -2.15%  '

Set-TextValue $ws.Range('B19') 'WrappedBTC'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D19') '69.671.70'
Set-TextValue $ws.Range('E19') '  -0.39%  '

Set-TextValue $ws.Range('E20') '  +0.93%  '

Set-TextValue $ws.Range('D21') '17.14'
Set-TextValue $ws.Range('E21') '  -0.15%  '

Set-TextValue $ws.Range('D22') '0.852'
Set-TextValue $ws.Range('E22') '  -1.77%  '

Set-TextValue $ws.Range('D23') '8.90'
Set-TextValue $ws.Range('E23') '  -3.30%  '

Set-TextValue $ws.Range('B24') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D24') '15.18'
Set-TextValue $ws.Range('E24') '  -2.28%  '

Set-TextValue $ws.Range('B25') 'Litecoin'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D25') '95.18'
Set-TextValue $ws.Range('E25') '  -0.48%  '

Set-TextValue $ws.Range('D26') '3.62'
Set-TextValue $ws.Range('E26') '  -2.15%  '

Set-TextValue $ws.Range('D27') '1.00'
Set-TextValue $ws.Range('E27') '  -0.03%  '

Set-TextValue $ws.Range('E28') '  -4.30%  '

Set-TextValue $ws.Range('D29') '32.87'
Set-TextValue $ws.Range('E29') '  -2.80%  '

Set-TextValue $ws.Range('D30') '8.64'
Set-TextValue $ws.Range('E30') '  -3.56%  '

Set-TextValue $ws.Range('D31') '7.84'
Set-TextValue $ws.Range('E31') '  -2.94%  '

Set-TextValue $ws.Range('D32') '2.82'
Set-TextValue $ws.Range('E32') '  -6.52%  '

Set-TextValue $ws.Range('E33') '  -3.21%  '

Set-TextValue $ws.Range('D34') '6.56'
Set-TextValue $ws.Range('E34') '  -5.22%  '

Set-TextValue $ws.Range('D35') '577.32'
Set-TextValue $ws.Range('E35') '  -19.77%  '

Set-TextValue $ws.Range('B36') 'VeChain'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D36') '0.0482'
Set-TextValue $ws.Range('E36') '  +2.01%  '

Set-TextValue $ws.Range('B37') 'Cosmos'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D37') '10.60'
Set-TextValue $ws.Range('E37') '  -0.50%  '

Set-TextValue $ws.Range('D38') '0.0959'
Set-TextValue $ws.Range('E38') '  -3.58%  '

Set-TextValue $ws.Range('D39') '1.00'
Set-TextValue $ws.Range('E39') '  +0.49%  '

Set-TextValue $ws.Range('E40') '  -0.79%  '

Set-TextValue $ws.Range('D41') '0.140'
Set-TextValue $ws.Range('E41') '  -1.30%  '

Set-TextValue $ws.Range('D42') '3.13'
Set-TextValue $ws.Range('E42') '  -11.41%  '

Set-TextValue $ws.Range('D43') '3.247.73'
Set-TextValue $ws.Range('E43') '  -2.94%  '

Set-TextValue $ws.Range('D44') '0.0₃0695'
Set-TextValue $ws.Range('E44') '  +1.12%  '

Set-TextValue $ws.Range('D45') '0.296'
Set-TextValue $ws.Range('E45') '  -5.80%  '

Set-TextValue $ws.Range('D46') '30.71'
Set-TextValue $ws.Range('E46') '  -4.81%  '

Set-TextValue $ws.Range('D47') '2.76'
Set-TextValue $ws.Range('E47') '  -4.53%  '

Set-TextValue $ws.Range('D48') '2.40'
Set-TextValue $ws.Range('E48') '  -5.95%  '

Set-TextValue $ws.Range('E49') '  -2.85%  '

Set-TextValue $ws.Range('D50') '133.35'
Set-TextValue $ws.Range('E50') '  +0.33%  '
